$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{A='AK_1.png'; B=0.576; C=0.002; D=0.981; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='AK_2.png'; B=0.883; C=0.003; D=0.993; E=0.571; F='Akhlak Kamiswara'; G='Benar'}
    @{A='AK_3.png'; B=0.859; C=0.003; D=0.985; E=0.857; F='Akhlak Kamiswara'; G='Benar'}
    @{A='AK_4.png'; B=0.536; C=0.002; D=0.986; E=0.571; F='Muhammad Iqbal Baqi'; G='Salah'}
    @{A='AK_5.png'; B=0.539; C=0.002; D=0.992; E=0.571; F='Akhlak Kamiswara'; G='Benar'}
    @{A='MIB_1.png'; B=1.601; C=0.005; D=0.979; E=0.857; F='Muhammad Iqbal Baqi'; G='Benar'}
    @{A='MIB_2.png'; B=1.165; C=0.004; D=0.982; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='MIB_3.png'; B=1.706; C=0.006; D=1; E=1; F='Muhammad Iqbal Baqi'; G='Benar'}
    @{A='MIB_4.png'; B=1.163; C=0.004; D=0.989; E=0.857; F='Muhammad Iqbal Baqi'; G='Benar'}
    @{A='MIB_5.png'; B=1.49; C=0.005; D=0.988; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='AAH_1.png'; B=0.872; C=0.003; D=0.984; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='AAH_2.png'; B=1.1; C=0.004; D=0.995; E=0.857; F='Andrea Ayunove Hutami'; G='Benar'}
    @{A='AAH_3.png'; B=0.909; C=0.003; D=0.979; E=0.714; F='Andrea Ayunove Hutami'; G='Benar'}
    @{A='TI_1.png'; B=0.988; C=0.003; D=0.966; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='TI_2.png'; B=0.957; C=0.003; D=0.994; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='TI_3.png'; B=0.758; C=0.003; D=0.994; E=1; F='Toni Ismail'; G='Benar'}
    @{A='TI_4.png'; B=0.737; C=0.002; D=0.99; E=1; F='Toni Ismail'; G='Benar'}
    @{A='TI_5.png'; B=1.043; C=0.003; D=0.996; E=0.571; F='Toni Ismail'; G='Benar'}
    @{A='RAS_1.png'; B=0.6820000000000001; C=0.002; D=0.984; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='RAS_2.png'; B=1.031; C=0.003; D=0.973; E=0.857; F='Ridha Ayu Salsabila'; G='Benar'}
    @{A='RAS_3.png'; B=0.5639999999999999; C=0.002; D=0.987; E=0.571; F='Muhammad Iqbal Baqi'; G='Salah'}
    @{A='RAS_4.png'; B=1.253; C=0.004; D=0.979; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='RAS_5.png'; B=1.182; C=0.004; D=0.98; E=0.857; F='Ridha Ayu Salsabila'; G='Benar'}
    @{A='RR_1.png'; B=1.437; C=0.005; D=0.985; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='RR_2.png'; B=1.381; C=0.005; D=0.99; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='RR_3.png'; B=1.167; C=0.004; D=0.92; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='RR_4.png'; B=1.373; C=0.005; D=0.991; E=0.571; F='Rafiqo Rapitasari'; G='Benar'}
    @{A='RR_5.png'; B=1.418; C=0.005; D=0.987; E=0.857; F='Rafiqo Rapitasari'; G='Benar'}
    @{A='AR_1.png'; B=0.886; C=0.003; D=0.976; E=1; F='Arizli Romadhon'; G='Benar'}
    @{A='GA_1.png'; B=1.547; C=0.005; D=0.986; E=1; F='Gege Ardiyansyah'; G='Benar'}
    @{A='GA_2.png'; B=0.674; C=0.002; D=0.99; E=0.571; F='Muhammad Iqbal Baqi'; G='Salah'}
    @{A='GA_3.png'; B=0.72; C=0.002; D=0.995; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='FY_1.png'; B=1.148; C=0.004; D=0.979; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='FY_2.png'; B=1.63; C=0.005; D=0.986; E=0.286; F='Tidak Diketahui'; G='Salah'}
    @{A='FY_3.png'; B=1.456; C=0.005; D=0.992; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='FY_4.png'; B=1.251; C=0.004; D=0.99; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='TO_1.png'; B=0.8070000000000001; C=0.003; D=0.978; E=0.286; F='Tidak Diketahui'; G='Salah'}
    @{A='TO_2.png'; B=1.072; C=0.004; D=0.988; E=0.857; F='Tiara Oktavian'; G='Benar'}
    @{A='TO_3.png'; B=0.836; C=0.003; D=0.992; E=0.714; F='Tiara Oktavian'; G='Benar'}
    @{A='TO_4.png'; B=3.738; C=0.012; D=0.857; E=0.429; F='Tidak Diketahui'; G='Salah'}
    @{A='TO_5.png'; B=3.231; C=0.011; D=0.849; E=0.571; F='Muhammad Iqbal Baqi'; G='Salah'}
    @{A='TD_1.png'; B=2.381; C=0.008; D=0.868; E=0.857; F='Muhammad Iqbal Baqi'; G='Salah'}
    @{A='TD_2.png'; B=2.561; C=0.008; D=0.869; E=0.429; F='Tidak Diketahui'; G='Benar'}
    @{A='TD_3.png'; B=1.051; C=0.004; D=0.972; E=0.429; F='Tidak Diketahui'; G='Benar'}
    @{A='TD_4.png'; B=1.003; C=0.003; D=0.991; E=0.286; F='Tidak Diketahui'; G='Benar'}
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
}
